# Scheduled runner refresh: updates Leve profit columns (H:N) with newly
# polled Universalis market prices for the affected Leve rows across all
# job sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet="ALC"; Row=41; Cells=@{"H"=182.78261; "I"=229.09091; "J"=140.33333; "K"=229.09091; "L"=140.33333; "M"=210.90909; "N"=-1020.33333} }
    @{ Sheet="ALC"; Row=62; Cells=@{"H"=166669330; "I"=500001500; "J"=3253; "K"=500001500; "L"=3253; "M"=-500000876; "N"=-4501} }
    @{ Sheet="ALC"; Row=64; Cells=@{"H"=3877.5; "I"=3633.3333; "J"=3982.1428; "K"=3633.3333; "L"=3982.1428; "M"=-3385.3333; "N"=-4478.1428} }
    @{ Sheet="ALC"; Row=65; Cells=@{"H"=166669330; "I"=500001500; "J"=3253; "K"=2500007500; "L"=16265; "M"=-2500004380; "N"=-22505} }
    @{ Sheet="ALC"; Row=67; Cells=@{"H"=3877.5; "I"=3633.3333; "J"=3982.1428; "K"=3633.3333; "L"=3982.1428; "M"=-2775.3333; "N"=-5698.1428} }
    @{ Sheet="ALC"; Row=86; Cells=@{"H"=1633.238; "I"=1653.1177; "J"=1548.75; "K"=1653.1177; "L"=1548.75; "M"=-530.1177; "N"=-3794.75} }
    @{ Sheet="ALC"; Row=89; Cells=@{"H"=1633.238; "I"=1653.1177; "J"=1548.75; "K"=8265.5885; "L"=7743.75; "M"=-2649.5885; "N"=-18975.75} }
    @{ Sheet="ALC"; Row=116; Cells=@{"H"=299753.22; "I"=4617; "J"=460736.6; "K"=4617; "L"=460736.6; "M"=-1175; "N"=-467620.6} }
    @{ Sheet="ALC"; Row=132; Cells=@{"H"=18502384; "I"=27751930; "J"=3290.9; "K"=83255790; "L"=9872.700000000001; "M"=-83253260; "N"=-14932.7} }
    @{ Sheet="ALC"; Row=134; Cells=@{"H"=75400; "J"=75400; "L"=75400; "N"=-85540} }
    @{ Sheet="ALC"; Row=135; Cells=@{"H"=1894.2727; "I"=0; "J"=1894.2727; "K"=0; "L"=17048.4543; "M"=$null; "N"=-22118.4543} }
    @{ Sheet="ALC"; Row=137; Cells=@{"H"=1319.88; "I"=1041.7858; "K"=3125.3574; "M"=-575.3574000000003} }
    @{ Sheet="ALC"; Row=141; Cells=@{"H"=9641; "I"=0; "K"=0; "M"=$null} }
    @{ Sheet="ARM"; Row=19; Cells=@{"H"=2000; "I"=2000; "K"=2000; "M"=-1771} }
    @{ Sheet="ARM"; Row=61; Cells=@{"H"=2847.1538; "I"=0; "J"=2847.1538; "K"=0; "L"=2847.1538; "M"=$null; "N"=-3271.1538} }
    @{ Sheet="ARM"; Row=74; Cells=@{"H"=896.925; "I"=559.56525; "J"=1353.3529; "K"=559.56525; "L"=1353.3529; "M"=314.43475; "N"=-3101.3529} }
    @{ Sheet="ARM"; Row=77; Cells=@{"H"=896.925; "I"=559.56525; "J"=1353.3529; "K"=2797.82625; "L"=6766.7645; "M"=1570.17375; "N"=-15502.7645} }
    @{ Sheet="ARM"; Row=92; Cells=@{"H"=19950; "J"=19950; "L"=19950; "N"=-24942} }
    @{ Sheet="ARM"; Row=132; Cells=@{"H"=1787286.9; "I"=2778664.8; "J"=2806.8; "K"=8335994.399999999; "L"=8420.400000000001; "M"=-8333464.399999999; "N"=-13480.4} }
    @{ Sheet="ARM"; Row=136; Cells=@{"H"=2847.1538; "I"=0; "J"=2847.1538; "K"=0; "L"=8541.4614; "M"=$null; "N"=-13641.4614} }
    @{ Sheet="BSM"; Row=22; Cells=@{"H"=768.36365; "I"=810.2; "J"=350; "K"=810.2; "L"=350; "M"=-637.2; "N"=-696} }
    @{ Sheet="BSM"; Row=94; Cells=@{"H"=982.2353000000001; "I"=589.8182; "J"=1701.6666; "K"=589.8182; "L"=1701.6666; "M"=-138.8182; "N"=-2603.6666} }
    @{ Sheet="BSM"; Row=134; Cells=@{"H"=9820183; "I"=15898209; "J"=1832.5385; "K"=47694627; "L"=5497.6155; "M"=-47692092; "N"=-10567.6155} }
    @{ Sheet="CRP"; Row=31; Cells=@{"H"=8377.758; "I"=1552.0526; "J"=17641.215; "K"=1552.0526; "L"=17641.215; "M"=-1257.0526; "N"=-18231.215} }
    @{ Sheet="CRP"; Row=34; Cells=@{"H"=8377.758; "I"=1552.0526; "J"=17641.215; "K"=1552.0526; "L"=17641.215; "M"=-1350.0526; "N"=-18045.215} }
    @{ Sheet="CRP"; Row=58; Cells=@{"H"=3429017.5; "I"=5534269.5; "J"=7983.4375; "K"=5534269.5; "L"=7983.4375; "M"=-5534066.5; "N"=-8389.4375} }
    @{ Sheet="CRP"; Row=132; Cells=@{"H"=11116333; "I"=30304196; "J"=7570.0527; "K"=90912588; "L"=22710.1581; "M"=-90910058; "N"=-27770.1581} }
    @{ Sheet="CRP"; Row=134; Cells=@{"H"=8014627.5; "I"=14707809; "J"=2842623.8; "K"=44123427; "L"=8527871.399999999; "M"=-44120892; "N"=-8532941.399999999} }
    @{ Sheet="CRP"; Row=136; Cells=@{"H"=3429017.5; "I"=5534269.5; "J"=7983.4375; "K"=16602808.5; "L"=23950.3125; "M"=-16600258.5; "N"=-29050.3125} }
    @{ Sheet="GSM"; Row=132; Cells=@{"H"=38505428; "I"=77002216; "J"=8643.23; "K"=231006648; "L"=25929.69; "M"=-231004118; "N"=-30989.69} }
    @{ Sheet="LTW"; Row=132; Cells=@{"H"=4446509.5; "I"=6251557.5; "K"=18754672.5; "M"=-18752142.5} }
    @{ Sheet="LTW"; Row=136; Cells=@{"H"=3474.721; "I"=3453.1765; "J"=3556.111; "K"=10359.5295; "L"=10668.333; "M"=-7809.529500000001; "N"=-15768.333} }
    @{ Sheet="WVR"; Row=86; Cells=@{"H"=10000.5; "J"=10000.5; "L"=10000.5; "N"=-12246.5} }
    @{ Sheet="WVR"; Row=89; Cells=@{"H"=10000.5; "J"=10000.5; "L"=50002.5; "N"=-61234.5} }
    @{ Sheet="WVR"; Row=132; Cells=@{"H"=46726290; "I"=0; "J"=46726290; "K"=0; "L"=140178870; "M"=$null; "N"=-140183930} }
    @{ Sheet="WVR"; Row=136; Cells=@{"H"=15037317; "I"=8304818; "J"=35715704; "K"=24914454; "L"=107147112; "M"=-24911904; "N"=-107152212} }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    foreach ($col in $u.Cells.Keys) {
        $val = $u.Cells[$col]
        $addr = "$col$($u.Row)"
        if ($null -eq $val) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}
